# Update "想去人数" (want-to-go count) values across sheets to reflect
# regenerated output data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1157
$ws.Range("F5").Value = 64
$ws.Range("F6").Value = 479
$ws.Range("F7").Value = 816
$ws.Range("F8").Value = 430
$ws.Range("F9").Value = 56
$ws.Range("F10").Value = 2060
$ws.Range("F12").Value = 239
$ws.Range("F14").Value = 958
$ws.Range("F15").Value = 135
$ws.Range("F16").Value = 2084
$ws.Range("F17").Value = 555
$ws.Range("F18").Value = 9588
$ws.Range("F19").Value = 935
$ws.Range("F25").Value = 147

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 6
$ws.Range("F10").Value = 134
$ws.Range("F13").Value = 39

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5636
$ws.Range("F3").Value = 445
$ws.Range("F4").Value = 417

# Sheet: 全部类型 (All Types - combined)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5636
$ws.Range("F4").Value = 445
$ws.Range("F5").Value = 417
$ws.Range("F7").Value = 1157
$ws.Range("F8").Value = 6
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 479
$ws.Range("F12").Value = 816
$ws.Range("F14").Value = 430
$ws.Range("F16").Value = 2060
$ws.Range("F18").Value = 239
$ws.Range("F22").Value = 958
$ws.Range("F24").Value = 135
$ws.Range("F25").Value = 134
$ws.Range("F27").Value = 2084
$ws.Range("F28").Value = 555
$ws.Range("F30").Value = 39
$ws.Range("F31").Value = 935
